# Applies the mob2dsl user_input.xlsx edit:
#  - DSL is now always produced with "ng" for human and "cp" for templates
#    (as opposed to "K cp"): every row's "Template" lookup column (F) on
#    the "allocation" sheet now points at the single Human DNA template
#    value instead of a per-row Promega lot number -> collapse the F
#    column to one value, which lets the extra shared strings get
#    garbage-collected on save.
#  - Input excel can have empty rows: an empty row is inserted after each
#    "set" of rows on the allocation sheet (before row 8 and before what
#    was row 14), matching the new blank-row-tolerant reader.
#  - Selection / active-sheet bookkeeping left by the editing session is
#    updated to match (allocation tab active w/ C9 selected, sample
#    layout tab shows B2 selected).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("allocation")
$ws2 = $wb.Worksheets.Item("sample layout")

# --- allocation sheet: make room for the two new blank rows -------------
# Row 8 becomes a blank separator row; everything that used to start at
# row 8 shifts down by one.
$ws1.Rows.Item(8).Insert()

# A second blank separator row is inserted right before the old row 14
# (now row 15, after the first insert above).
$ws1.Rows.Item(15).Insert()

# --- allocation sheet: collapse the per-row template lookup -------------
# Every populated row in column F ("Template" values, shared strings
# 33-45) now uses the single "HgDNA_Promega305466" value. Rows 8 and 15
# are the freshly-inserted blank separators and stay empty, as does the
# pre-existing blank separator (old row 7).
$templateValue = $ws1.Cells.Item(2, 6).Value2

$dataRows = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14, 16, 17)
foreach ($r in $dataRows) {
    $ws1.Cells.Item($r, 6).Value2 = $templateValue
}

# --- selection / active sheet bookkeeping --------------------------------
$ws2.Activate()
$ws2.Range("B2").Select()

$ws1.Activate()
$ws1.Range("C9").Select()
